$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$lo = $ws.ListObjects.Item(1)

# Insert a new column before the current column E (en_variable-label),
# shifting en_variable-label -> F and en_note -> G. Data (formatting/values)
# in the inserted column starts out blank.
$ws.Columns.Item(5).Insert()

# Grow the table definition to include the freshly inserted column.
$lo.Resize($ws.Range("A1:G4"))

# Re-assert the header text for every column so the table's column-name
# cache is resynced against the sheet (the newly inserted column becomes
# "timeval"; the two columns that shifted right need their names restated).
$ws.Range("E1").Value = "timeval"
$ws.Range("F1").Value = "en_variable-label"
$ws.Range("G1").Value = "en_note"

# variable-type (column D) no longer applies to the gender/time rows; the
# new timeval column takes over marking which variable is the time variable.
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()

$ws.Range("E2").Value = $true
$ws.Range("E3").Value = $true
$ws.Range("E4").Value = $false

# Column widths: D and E share a generic (non-autofit) width, and the
# en_variable-label column (now F) keeps its old best-fit width.
$ws.Columns.Item(4).ColumnWidth = 10.8
$ws.Columns.Item(5).ColumnWidth = 10.8

$ws.Range("D2").Select() | Out-Null
